$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-14 21:18:50'
$ws.Range('I2').Value = '35.5 mm'
$ws.Range('O2').Value = '-1.2 °C'
$ws.Range('E3').Value = '2026-02-14 21:18:52'
$ws.Range('I3').Value = '15.1 mm'
$ws.Range('L3').Value = '57.2 km/h - 249º 20:58 TU'
$ws.Range('N3').Value = '-7.3 °C 20:33 TU'
$ws.Range('E4').Value = '2026-02-14 21:18:55'
$ws.Range('H4').Value = "'72%"
$ws.Range('J4').Value = '997.4 hPa'
$ws.Range('E5').Value = '2026-02-14 21:18:57'
$ws.Range('I5').Value = '22.1 mm'
$ws.Range('N5').Value = '-7.0 °C 20:59 TU'
$ws.Range('E6').Value = '2026-02-14 21:19:00'
$ws.Range('H6').Value = "'74%"
$ws.Range('J6').Value = '997.5 hPa'
$ws.Range('O6').Value = '10.3 °C'
$ws.Range('E7').Value = '2026-02-14 21:19:02'
$ws.Range('H7').Value = "'51%"
$ws.Range('J7').Value = '997.6 hPa'
$ws.Range('E8').Value = '2026-02-14 21:19:05'
$ws.Range('J8').Value = '997.5 hPa'
$ws.Range('E9').Value = '2026-02-14 21:19:08'
$ws.Range('N9').Value = '10.0 °C 20:58 TU'
$ws.Range('E10').Value = '2026-02-14 21:19:11'
$ws.Range('H10').Value = "'76%"
$ws.Range('E11').Value = '2026-02-14 21:19:14'
$ws.Range('O11').Value = '6.3 °C'
$ws.Range('E12').Value = '2026-02-14 21:19:16'
$ws.Range('N12').Value = '10.4 °C 20:38 TU'
$ws.Range('O12').Value = '12.1 °C'
$ws.Range('E13').Value = '2026-02-14 21:19:18'
$ws.Range('J13').Value = '1000.3 hPa'
$ws.Range('O13').Value = '3.9 °C'
$ws.Range('E14').Value = '2026-02-14 21:19:21'
$ws.Range('O14').Value = '13.5 °C'
$ws.Range('E15').Value = '2026-02-14 21:19:23'
$ws.Range('N15').Value = '9.5 °C 20:59 TU'
$ws.Range('E16').Value = '2026-02-14 21:19:26'
$ws.Range('N16').Value = '-9.0 °C 20:42 TU'
$ws.Range('E17').Value = '2026-02-14 21:19:28'
$ws.Range('H17').Value = "'67%"
$ws.Range('N17').Value = '-0.9 °C 20:56 TU'
$ws.Range('O17').Value = '1.7 °C'
$ws.Range('E18').Value = '2026-02-14 21:19:31'
$ws.Range('H18').Value = "'75%"
$ws.Range('J18').Value = '997.7 hPa'
$ws.Range('L18').Value = '27.0 km/h - 13º 20:46 TU'
$ws.Range('E19').Value = '2026-02-14 21:19:34'
$ws.Range('H19').Value = "'76%"
$ws.Range('E20').Value = '2026-02-14 21:19:37'
$ws.Range('I20').Value = '4.3 mm'
$ws.Range('N20').Value = '-8.0 °C 20:59 TU'
$ws.Range('O20').Value = '-5.4 °C'
$ws.Range('E21').Value = '2026-02-14 21:19:40'
$ws.Range('J21').Value = '1000.2 hPa'
$ws.Range('L21').Value = '65.5 km/h - 6º 20:39 TU'
$ws.Range('E22').Value = '2026-02-14 21:19:43'
$ws.Range('I22').Value = '0.9 mm'
$ws.Range('N22').Value = '-9.5 °C 20:41 TU'
$ws.Range('E23').Value = '2026-02-14 21:19:46'
$ws.Range('I23').Value = '39.9 mm'
$ws.Range('N23').Value = '-8.4 °C 20:55 TU'
$ws.Range('O23').Value = '-6.1 °C'
$ws.Range('E24').Value = '2026-02-14 21:19:48'
$ws.Range('J24').Value = '1001.7 hPa'
$ws.Range('E25').Value = '2026-02-14 21:19:51'
$ws.Range('I25').Value = '16.4 mm'
$ws.Range('N25').Value = '-7.3 °C 20:57 TU'
$ws.Range('O25').Value = '-4.8 °C'
$ws.Range('E26').Value = '2026-02-14 21:19:54'
$ws.Range('E27').Value = '2026-02-14 21:19:56'
$ws.Range('O27').Value = '-3.2 °C'
$ws.Range('E28').Value = '2026-02-14 21:19:59'
$ws.Range('H28').Value = "'66%"
$ws.Range('J28').Value = '997.4 hPa'
$ws.Range('O28').Value = '9.3 °C'
$ws.Range('E29').Value = '2026-02-14 21:20:02'
$ws.Range('H29').Value = "'62%"
$ws.Range('E30').Value = '2026-02-14 21:20:05'
$ws.Range('J30').Value = '997.3 hPa'
$ws.Range('E31').Value = '2026-02-14 21:20:08'
$ws.Range('H31').Value = "'68%"
$ws.Range('J31').Value = '996.4 hPa'
$ws.Range('N31').Value = '8.0 °C 20:59 TU'
$ws.Range('O31').Value = '9.2 °C'
$ws.Range('E32').Value = '2026-02-14 21:20:10'
$ws.Range('E33').Value = '2026-02-14 21:20:13'
$ws.Range('J33').Value = '999.7 hPa'
$ws.Range('O33').Value = '3.9 °C'
$ws.Range('E34').Value = '2026-02-14 21:20:16'
$ws.Range('I34').Value = '3.7 mm'
$ws.Range('N34').Value = '-5.2 °C 20:55 TU'
$ws.Range('O34').Value = '-2.4 °C'
$ws.Range('E35').Value = '2026-02-14 21:20:18'
$ws.Range('J35').Value = '1004.3 hPa'
$ws.Range('N35').Value = '1.4 °C 20:56 TU'
$ws.Range('E36').Value = '2026-02-14 21:20:21'
$ws.Range('J36').Value = '998.1 hPa'
$ws.Range('N36').Value = '10.3 °C 20:57 TU'
$ws.Range('E37').Value = '2026-02-14 21:20:24'
$ws.Range('H37').Value = "'64%"
$ws.Range('J37').Value = '998.4 hPa'
$ws.Range('E38').Value = '2026-02-14 21:20:27'
$ws.Range('H38').Value = "'81%"
$ws.Range('E39').Value = '2026-02-14 21:20:30'
$ws.Range('I39').Value = '13.5 mm'
$ws.Range('N39').Value = '-8.5 °C 20:53 TU'
$ws.Range('E40').Value = '2026-02-14 21:20:32'
$ws.Range('J40').Value = '1000.9 hPa'
$ws.Range('E41').Value = '2026-02-14 21:20:35'
$ws.Range('J41').Value = '999.4 hPa'
$ws.Range('E42').Value = '2026-02-14 21:20:38'
$ws.Range('E43').Value = '2026-02-14 21:20:41'
$ws.Range('H43').Value = "'65%"
$ws.Range('E44').Value = '2026-02-14 21:20:43'
$ws.Range('G44').Value = '274 cm'
$ws.Range('I44').Value = '37.6 mm'
$ws.Range('N44').Value = '-7.6 °C 20:59 TU'
$ws.Range('E45').Value = '2026-02-14 21:20:46'
$ws.Range('H45').Value = "'83%"
$ws.Range('I45').Value = '13.8 mm'
$ws.Range('J45').Value = '1006.7 hPa'
$ws.Range('N45').Value = '0.5 °C 20:52 TU'
$ws.Range('O45').Value = '2.9 °C'
$ws.Range('E46').Value = '2026-02-14 21:20:48'
$ws.Range('J46').Value = '1002.9 hPa'
$ws.Range('O46').Value = '11.6 °C'
